$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Header paragraph: "suscrito con fecha  24 de mayo de 2022" -> "... 31 ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "suscrito con fecha  24 de mayo de 2022",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "suscrito con fecha  31 de mayo de 2022", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Final signature paragraph: "En Puertollano a  24  de mayo  2022" -> "... 31 ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "En Puertollano a  24  de mayo  2022",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "En Puertollano a  31  de mayo  2022", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Students data table (second table in the document): add grey cell
#    borders to every cell, and update the cell contents.
# ---------------------------------------------------------------------------
$wdBorderTop    = -1
$wdBorderLeft   = -2
$wdBorderBottom = -3
$wdBorderRight  = -4
$wdLineStyleSingle = 1
$borderWidth = 3          # -> w:sz="6" (single, 0.75pt) in the saved OOXML
$borderColor = 8947848    # 0x888888

$table = $d.Tables.Item(2)

function Set-CellBorders($cell) {
    $cell.Borders.Item($wdBorderTop).LineStyle = $wdLineStyleSingle
    $cell.Borders.Item($wdBorderTop).LineWidth = $borderWidth
    $cell.Borders.Item($wdBorderTop).Color = $borderColor

    $cell.Borders.Item($wdBorderLeft).LineStyle = $wdLineStyleSingle
    $cell.Borders.Item($wdBorderLeft).LineWidth = $borderWidth
    $cell.Borders.Item($wdBorderLeft).Color = $borderColor

    $cell.Borders.Item($wdBorderRight).LineStyle = $wdLineStyleSingle
    $cell.Borders.Item($wdBorderRight).LineWidth = $borderWidth
    $cell.Borders.Item($wdBorderRight).Color = $borderColor

    $cell.Borders.Item($wdBorderBottom).LineStyle = $wdLineStyleSingle
    $cell.Borders.Item($wdBorderBottom).LineWidth = $borderWidth
    $cell.Borders.Item($wdBorderBottom).Color = $borderColor
}

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        Set-CellBorders $table.Cell($r, $c)
    }
}

# Row 2 (first data row) content: Moreno Ramos Laura / 13c / dsfds / 2022-05-24 / 2022-05-26
#   -> Díez Viñas Malena / 14d / dsf / 2022-05-31 / 2022-06-01
$table.Cell(2, 1).Range.Text = "Díez Viñas Malena"
$table.Cell(2, 2).Range.Text = "14d"
$table.Cell(2, 4).Range.Text = "dsf"
$table.Cell(2, 6).Range.Text = "2022-05-31"
$table.Cell(2, 7).Range.Text = "2022-06-01"

# Row 3 (second data row) content: Díez Viñas Malena / 14d / sdfdsf / 2022-05-24 / 2022-05-26
#   -> Moreno Ramos Laura / 13c / dsfdsf / 2022-05-31 / 2022-06-01
$table.Cell(3, 1).Range.Text = "Moreno Ramos Laura"
$table.Cell(3, 2).Range.Text = "13c"
$table.Cell(3, 4).Range.Text = "dsfdsf"
$table.Cell(3, 6).Range.Text = "2022-05-31"
$table.Cell(3, 7).Range.Text = "2022-06-01"

Write-Output "done"
